# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 53 ("Región de O'Higgins", 2023-04-05),
# pushing the previous row 53 ("Provincia de Limarí", 2021-04-30) down to row 54.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 53, shifting the former row 53 (and below) down.
$ws.Rows.Item(53).Insert()

# Populate the newly inserted row 53 with the new weekly data.
$ws.Cells.Item(53, 1).Value = 10
$ws.Cells.Item(53, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(53, 3).Value = "La Araucanía"
$ws.Cells.Item(53, 4).Value = 45021
$ws.Cells.Item(53, 5).Value = 9
$ws.Cells.Item(53, 6).Value = "Fruta"
$ws.Cells.Item(53, 7).Value = 100107
$ws.Cells.Item(53, 8).Value = "Otros"
$ws.Cells.Item(53, 9).Value = 100107001
$ws.Cells.Item(53, 10).Value = "Caqui"
$ws.Cells.Item(53, 11).Value = "Fuyu"
$ws.Cells.Item(53, 12).Value = "Primera"
$ws.Cells.Item(53, 13).Value = 90
$ws.Cells.Item(53, 14).Value = 24000
$ws.Cells.Item(53, 15).Value = 24000
$ws.Cells.Item(53, 16).Value = 24000
$ws.Cells.Item(53, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(53, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(53, 19).Value = 1600
$ws.Cells.Item(53, 20).Value = 15
